$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the lone cell in row 9 (it becomes a fully empty row)
$ws.Range("D9").Clear()

# Remove the cells D10:F10 entirely (contribution scores no longer entered)
$ws.Range("D10:F10").Clear()

# Clear the remaining contribution-score values in row 10 but keep their styling
$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("G10").ClearContents()

# Move the active selection to B10
[void]$ws.Range("B10").Select()
